$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 702.34784
$ws.Range("I28").Value = 598.1053000000001
$ws.Range("J28").Value = 1197.5
$ws.Range("K28").Value = 598.1053000000001
$ws.Range("L28").Value = 1197.5
$ws.Range("M28").Value = -113.1053000000001
$ws.Range("N28").Value = -2167.5

$ws.Range("H62").Value = 23989.1
$ws.Range("I62").Value = 35448.5
$ws.Range("J62").Value = 6800
$ws.Range("K62").Value = 35448.5
$ws.Range("L62").Value = 6800
$ws.Range("M62").Value = -34824.5
$ws.Range("N62").Value = -8048

$ws.Range("H65").Value = 23989.1
$ws.Range("I65").Value = 35448.5
$ws.Range("J65").Value = 6800
$ws.Range("K65").Value = 177242.5
$ws.Range("L65").Value = 34000
$ws.Range("M65").Value = -174122.5
$ws.Range("N65").Value = -40240

$ws.Range("H107").Value = 558.7143
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 400
$ws.Range("M107").Value = 1520

$ws.Range("H137").Value = 3161.4482
$ws.Range("I137").Value = 3725
$ws.Range("J137").Value = 1001.1667
$ws.Range("K137").Value = 11175
$ws.Range("L137").Value = 3003.5001
$ws.Range("M137").Value = -8625
$ws.Range("N137").Value = -8103.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1454.2142
$ws.Range("I2").Value = 1556.174
$ws.Range("J2").Value = 985.2
$ws.Range("K2").Value = 1556.174
$ws.Range("L2").Value = 985.2
$ws.Range("M2").Value = -1443.174
$ws.Range("N2").Value = -1211.2

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H32").Value = 2688.62
$ws.Range("I32").Value = 1978.5333
$ws.Range("J32").Value = 9079.4
$ws.Range("K32").Value = 1978.5333
$ws.Range("L32").Value = 9079.4
$ws.Range("M32").Value = -1691.5333
$ws.Range("N32").Value = -9653.4

$ws.Range("H45").Value = 1311.7567
$ws.Range("I45").Value = 1153.619
$ws.Range("J45").Value = 1519.3125
$ws.Range("K45").Value = 1153.619
$ws.Range("L45").Value = 1519.3125
$ws.Range("M45").Value = -776.6189999999999
$ws.Range("N45").Value = -2273.3125

$ws.Range("H74").Value = 158762.47
$ws.Range("I74").Value = 218230.56
$ws.Range("J74").Value = 44781.957
$ws.Range("K74").Value = 218230.56
$ws.Range("L74").Value = 44781.957
$ws.Range("M74").Value = -217356.56
$ws.Range("N74").Value = -46529.957

$ws.Range("H77").Value = 158762.47
$ws.Range("I77").Value = 218230.56
$ws.Range("J77").Value = 44781.957
$ws.Range("K77").Value = 1091152.8
$ws.Range("L77").Value = 223909.785
$ws.Range("M77").Value = -1086784.8
$ws.Range("N77").Value = -232645.785

$ws.Range("H116").Value = 1454.2142
$ws.Range("I116").Value = 1556.174
$ws.Range("J116").Value = 985.2
$ws.Range("K116").Value = 1556.174
$ws.Range("L116").Value = 985.2
$ws.Range("M116").Value = 737.826
$ws.Range("N116").Value = -5573.2

$ws.Range("H122").Value = 2957.158
$ws.Range("I122").Value = 2499.1765
$ws.Range("J122").Value = 6850
$ws.Range("K122").Value = 7497.529500000001
$ws.Range("L122").Value = 20550
$ws.Range("M122").Value = -5047.529500000001
$ws.Range("N122").Value = -25450

$ws.Range("H132").Value = 1976.1428
$ws.Range("I132").Value = 1766.5454
$ws.Range("J132").Value = 2408.4375
$ws.Range("K132").Value = 5299.6362
$ws.Range("L132").Value = 7225.3125
$ws.Range("M132").Value = -2769.6362
$ws.Range("N132").Value = -12285.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1454.2142
$ws.Range("I3").Value = 1556.174
$ws.Range("J3").Value = 985.2
$ws.Range("K3").Value = 1556.174
$ws.Range("L3").Value = 985.2
$ws.Range("M3").Value = -1442.174
$ws.Range("N3").Value = -1213.2

$ws.Range("H20").Value = 2011.1111
$ws.Range("I20").Value = 2184.5
$ws.Range("J20").Value = 1924.4166
$ws.Range("K20").Value = 2184.5
$ws.Range("L20").Value = 1924.4166
$ws.Range("M20").Value = -1937.5
$ws.Range("N20").Value = -2418.4166

$ws.Range("H86").Value = 14148.25
$ws.Range("I86").Value = 10247.667
$ws.Range("J86").Value = 25850
$ws.Range("K86").Value = 10247.667
$ws.Range("L86").Value = 25850
$ws.Range("M86").Value = -9124.666999999999
$ws.Range("N86").Value = -28096

$ws.Range("H89").Value = 14148.25
$ws.Range("I89").Value = 10247.667
$ws.Range("J89").Value = 25850
$ws.Range("K89").Value = 51238.335
$ws.Range("L89").Value = 129250
$ws.Range("M89").Value = -45622.335
$ws.Range("N89").Value = -140482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2954
$ws.Range("I58").Value = 3193.756
$ws.Range("J58").Value = 2407.889
$ws.Range("K58").Value = 3193.756
$ws.Range("L58").Value = 2407.889
$ws.Range("M58").Value = -2990.756
$ws.Range("N58").Value = -2813.889

$ws.Range("H136").Value = 2954
$ws.Range("I136").Value = 3193.756
$ws.Range("J136").Value = 2407.889
$ws.Range("K136").Value = 9581.268
$ws.Range("L136").Value = 7223.667
$ws.Range("M136").Value = -7031.268
$ws.Range("N136").Value = -12323.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1920
$ws.Range("I137").Value = 1996.6666
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 5989.9998
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -889.9997999999996
$ws.Range("N137").Value = -13200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4099.643
$ws.Range("I132").Value = 4412.5264
$ws.Range("J132").Value = 3841.1738
$ws.Range("K132").Value = 13237.5792
$ws.Range("L132").Value = 11523.5214
$ws.Range("M132").Value = -10707.5792
$ws.Range("N132").Value = -16583.5214

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1707.8948
$ws.Range("I132").Value = 1181.9565
$ws.Range("J132").Value = 2201.6326
$ws.Range("K132").Value = 3545.8695
$ws.Range("L132").Value = 6604.8978
$ws.Range("M132").Value = -1015.8695
$ws.Range("N132").Value = -11664.8978

$ws.Range("H141").Value = 43905
$ws.Range("J141").Value = 43905
$ws.Range("L141").Value = 43905
$ws.Range("N141").Value = -54265
